$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 3097.9285
$ws.Range("I98").Value = 2235.3235
$ws.Range("J98").Value = 6764
$ws.Range("K98").Value = 2235.3235
$ws.Range("L98").Value = 6764
$ws.Range("M98").Value = -737.3235
$ws.Range("N98").Value = -9760
$ws.Range("H100").Value = 3119.4375
$ws.Range("I100").Value = 2690.5
$ws.Range("K100").Value = 2690.5
$ws.Range("M100").Value = -2149.5
$ws.Range("H106").Value = 3582
$ws.Range("I106").Value = 5955
$ws.Range("K106").Value = 5955
$ws.Range("M106").Value = -5324
$ws.Range("H122").Value = 3097.9285
$ws.Range("I122").Value = 2235.3235
$ws.Range("J122").Value = 6764
$ws.Range("K122").Value = 6705.970499999999
$ws.Range("L122").Value = 20292
$ws.Range("M122").Value = -4255.970499999999
$ws.Range("N122").Value = -25192
$ws.Range("H141").Value = 4016.35
$ws.Range("I141").Value = 1481.6842
$ws.Range("J141").Value = 52175
$ws.Range("K141").Value = 4445.0526
$ws.Range("L141").Value = 156525
$ws.Range("M141").Value = 734.9474
$ws.Range("N141").Value = -166885

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 20000
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 20000
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 20000
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value = -20288
$ws.Range("H102").Value = 66069.06
$ws.Range("I102").Value = 3500.3845
$ws.Range("J102").Value = 337200
$ws.Range("K102").Value = 3500.3845
$ws.Range("L102").Value = 337200
$ws.Range("M102").Value = -1878.3845
$ws.Range("N102").Value = -340444
$ws.Range("H122").Value = 4407.4136
$ws.Range("I122").Value = 4546.1177
$ws.Range("K122").Value = 13638.3531
$ws.Range("M122").Value = -11188.3531

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3200
$ws.Range("I105").Value = 3502
$ws.Range("J105").Value = 2822.5
$ws.Range("K105").Value = 3502
$ws.Range("L105").Value = 2822.5
$ws.Range("M105").Value = -1755
$ws.Range("N105").Value = -6316.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2059.611
$ws.Range("I31").Value = 1301.3914
$ws.Range("J31").Value = 3401.077
$ws.Range("K31").Value = 1301.3914
$ws.Range("L31").Value = 3401.077
$ws.Range("M31").Value = -1006.3914
$ws.Range("N31").Value = -3991.077
$ws.Range("H34").Value = 2059.611
$ws.Range("I34").Value = 1301.3914
$ws.Range("J34").Value = 3401.077
$ws.Range("K34").Value = 1301.3914
$ws.Range("L34").Value = 3401.077
$ws.Range("M34").Value = -1099.3914
$ws.Range("N34").Value = -3805.077
$ws.Range("H99").Value = 3030.1428
$ws.Range("I99").Value = 2951.8333
$ws.Range("J99").Value = 3500
$ws.Range("K99").Value = 2951.8333
$ws.Range("L99").Value = 3500
$ws.Range("M99").Value = -1453.8333
$ws.Range("N99").Value = -6496
$ws.Range("H105").Value = 2000
$ws.Range("I105").Value = 2000
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2000
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -253
$ws.Range("N105").ClearContents()
$ws.Range("H126").Value = 3030.1428
$ws.Range("I126").Value = 2951.8333
$ws.Range("J126").Value = 3500
$ws.Range("K126").Value = 8855.499899999999
$ws.Range("L126").Value = 10500
$ws.Range("M126").Value = -6385.499899999999
$ws.Range("N126").Value = -15440

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1372.6145
$ws.Range("I131").Value = 372.22223
$ws.Range("J131").Value = 1494.2838
$ws.Range("K131").Value = 1116.66669
$ws.Range("L131").Value = 4482.8514
$ws.Range("M131").Value = 3923.33331
$ws.Range("N131").Value = -14562.8514
$ws.Range("H132").Value = 1961
$ws.Range("I132").Value = 1562.6471
$ws.Range("J132").Value = 2576.6365
$ws.Range("K132").Value = 14063.8239
$ws.Range("L132").Value = 23189.7285
$ws.Range("M132").Value = -11533.8239
$ws.Range("N132").Value = -28249.7285
$ws.Range("H133").Value = 4347.227
$ws.Range("I133").Value = 1260.125
$ws.Range("J133").Value = 6111.2856
$ws.Range("K133").Value = 3780.375
$ws.Range("L133").Value = 18333.8568
$ws.Range("M133").Value = 1279.625
$ws.Range("N133").Value = -28453.8568

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 100000000
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H8").Value = 100000000
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()
$ws.Range("H97").Value = 163180
$ws.Range("I97").Value = 124840
$ws.Range("J97").Value = 201520
$ws.Range("K97").Value = 124840
$ws.Range("L97").Value = 201520
$ws.Range("M97").Value = -124344
$ws.Range("N97").Value = -202512
$ws.Range("H122").Value = 2582.3076
$ws.Range("I122").Value = 1956.3334
$ws.Range("J122").Value = 3990.75
$ws.Range("K122").Value = 5869.0002
$ws.Range("L122").Value = 11972.25
$ws.Range("M122").Value = -3419.0002
$ws.Range("N122").Value = -16872.25

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("N60").ClearContents()
$ws.Range("H69").Value = 103809
$ws.Range("I69").Value = 30000
$ws.Range("J69").Value = 116110.5
$ws.Range("K69").Value = 30000
$ws.Range("L69").Value = 116110.5
$ws.Range("M69").Value = -29189
$ws.Range("N69").Value = -117732.5
$ws.Range("H72").Value = 103809
$ws.Range("I72").Value = 30000
$ws.Range("J72").Value = 116110.5
$ws.Range("K72").Value = 90000
$ws.Range("L72").Value = 348331.5
$ws.Range("M72").Value = -85944
$ws.Range("N72").Value = -356443.5
$ws.Range("H93").Value = 1799.5
$ws.Range("I93").Value = 968.38464
$ws.Range("J93").Value = 3000
$ws.Range("K93").Value = 968.38464
$ws.Range("L93").Value = 3000
$ws.Range("M93").Value = 279.61536
$ws.Range("N93").Value = -5496
$ws.Range("H122").Value = 20459198
$ws.Range("I122").Value = 20836814
$ws.Range("K122").Value = 62510442
$ws.Range("M122").Value = -62507992

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1903.6875
$ws.Range("J96").Value = 2322.2222
$ws.Range("L96").Value = 2322.2222
$ws.Range("N96").Value = -5068.2222
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
